# Weekly update to the Jengibre (Vega Central Mapocho de Santiago) sheet:
# insert two new rows of the latest week's price data above the existing
# historical rows (rows 32-51 shift down to rows 34-53).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 32, pushing existing data down.
$ws.Rows("32:33").Insert()

# --- New row 32 (Primera) ---
$ws.Cells.Item(32, 1).Value2 = 9
$ws.Cells.Item(32, 2).Value2 = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(32, 3).Value2 = 'Metropolitana'
$ws.Cells.Item(32, 4).Value2 = 44438
$ws.Cells.Item(32, 5).Value2 = 13
$ws.Cells.Item(32, 6).Value2 = 100114007
$ws.Cells.Item(32, 7).Value2 = 'Jengibre'
$ws.Cells.Item(32, 8).Value2 = 'Sin especificar'
$ws.Cells.Item(32, 9).Value2 = 'Primera'
$ws.Cells.Item(32, 10).Value2 = 790
$ws.Cells.Item(32, 11).Value2 = 13000
$ws.Cells.Item(32, 12).Value2 = 14000
$ws.Cells.Item(32, 13).Value2 = 13494
$ws.Cells.Item(32, 14).Value2 = '$/caja 13 kilos'
$ws.Cells.Item(32, 15).Value2 = 'Perú'
$ws.Cells.Item(32, 16).Value2 = 1038
$ws.Cells.Item(32, 17).Value2 = 13
$ws.Cells.Item(32, 18).Value2 = 'Hortaliza'

# --- New row 33 (Segunda) ---
$ws.Cells.Item(33, 1).Value2 = 9
$ws.Cells.Item(33, 2).Value2 = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(33, 3).Value2 = 'Metropolitana'
$ws.Cells.Item(33, 4).Value2 = 44438
$ws.Cells.Item(33, 5).Value2 = 13
$ws.Cells.Item(33, 6).Value2 = 100114007
$ws.Cells.Item(33, 7).Value2 = 'Jengibre'
$ws.Cells.Item(33, 8).Value2 = 'Sin especificar'
$ws.Cells.Item(33, 9).Value2 = 'Segunda'
$ws.Cells.Item(33, 10).Value2 = 340
$ws.Cells.Item(33, 11).Value2 = 11000
$ws.Cells.Item(33, 12).Value2 = 12000
$ws.Cells.Item(33, 13).Value2 = 11500
$ws.Cells.Item(33, 14).Value2 = '$/caja 13 kilos'
$ws.Cells.Item(33, 15).Value2 = 'Perú'
$ws.Cells.Item(33, 16).Value2 = 885
$ws.Cells.Item(33, 17).Value2 = 13
$ws.Cells.Item(33, 18).Value2 = 'Hortaliza'

Write-Output "Inserted 2 new rows with latest week data; dimension now $($ws.UsedRange.Rows.Count) rows"
